$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 8, pushing "Jishi_explore_MySelf"/"End" rows down
$ws.Rows.Item(8).Insert()

# Fill in the newly inserted row 8 with the new test case data
$ws.Range("A8").Value = "List_banner"
$ws.Range("B8").Value = "//div[not(contains(@class,'ng-trigger ng-trigger-autoHeight ng-star-inserted cloned'))]/img[@class='img-fluid ng-star-inserted']"

# Widen column B to fit the new, longer content (best-fit width for the new text)
$ws.Columns.Item(2).ColumnWidth = 110.8

# Update the active selection to match the edited workbook
$ws.Range("B9").Select()
